{"js": "const replacements = [\n  { find: \"25\u00f75=\", replace: \"69\u00f79=\" },\n  { find: \"72\u00f76=\", replace: \"91\u00f72=\" },\n  { find: \"38\u00f78=\", replace: \"61\u00f77=\" },\n  { find: \"98\u00f76=\", replace: \"68\u00f79=\" },\n  { find: \"61\u00f74=\", replace: \"53\u00f75=\" },\n  { find: \"74\u00f74=\", replace: \"71\u00f77=\" },\n  { find: \"73\u00f79=\", replace: \"32\u00f75=\" },\n  { find: \"28\u00f72=\", replace: \"64\u00f77=\" },\n  { find: \"62\u00f76=\", replace: \"73\u00f76=\" },\n  { find: \"35\u00f76=\", replace: \"20\u00f77=\" },\n  { find: \"60\u00f76=\", replace: \"27\u00f76=\" },\n  { find: \"79\u00f74=\", replace: \"36\u00f76=\" },\n  { find: \"54\u00f72=\", replace: \"42\u00f77=\" },\n  { find: \"71\u00f74=\", replace: \"66\u00f79=\" },\n  { find: \"80\u00f79=\", replace: \"51\u00f77=\" },\n  { find: \"81\u00f79=\", replace: \"14\u00f74=\" },\n  { find: \"85\u00f72=\", replace: \"26\u00f74=\" },\n  { find: \"50\u00f79=\", replace: \"44\u00f74=\" },\n  { find: \"24\u00f73=\", replace: \"76\u00f73=\" },\n  { find: \"32\u00f76=\", replace: \"26\u00f79=\" },\n  { find: \"85\u00f79=\", replace: \"43\u00f79=\" },\n  { find: \"61\u00f73=\", replace: \"46\u00f77=\" },\n  { find: \"20\u00f75=\", replace: \"62\u00f72=\" },\n  { find: \"84\u00f76=\", replace: \"38\u00f74=\" },\n  { find: \"83\u00f75=\", replace: \"48\u00f74=\" },\n];\n\nconst body = context.document.body;\n\nfor (const { find, replace } of replacements) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replace, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"25\u00f75=\"; New = \"69\u00f79=\" },\n    @{ Old = \"72\u00f76=\"; New = \"91\u00f72=\" },\n    @{ Old = \"38\u00f78=\"; New = \"61\u00f77=\" },\n    @{ Old = \"98\u00f76=\"; New = \"68\u00f79=\" },\n    @{ Old = \"61\u00f74=\"; New = \"53\u00f75=\" },\n    @{ Old = \"74\u00f74=\"; New = \"71\u00f77=\" },\n    @{ Old = \"73\u00f79=\"; New = \"32\u00f75=\" },\n    @{ Old = \"28\u00f72=\"; New = \"64\u00f77=\" },\n    @{ Old = \"62\u00f76=\"; New = \"73\u00f76=\" },\n    @{ Old = \"35\u00f76=\"; New = \"20\u00f77=\" },\n    @{ Old = \"60\u00f76=\"; New = \"27\u00f76=\" },\n    @{ Old = \"79\u00f74=\"; New = \"36\u00f76=\" },\n    @{ Old = \"54\u00f72=\"; New = \"42\u00f77=\" },\n    @{ Old = \"71\u00f74=\"; New = \"66\u00f79=\" },\n    @{ Old = \"80\u00f79=\"; New = \"51\u00f77=\" },\n    @{ Old = \"81\u00f79=\"; New = \"14\u00f74=\" },\n    @{ Old = \"85\u00f72=\"; New = \"26\u00f74=\" },\n    @{ Old = \"50\u00f79=\"; New = \"44\u00f74=\" },\n    @{ Old = \"24\u00f73=\"; New = \"76\u00f73=\" },\n    @{ Old = \"32\u00f76=\"; New = \"26\u00f79=\" },\n    @{ Old = \"85\u00f79=\"; New = \"43\u00f79=\" },\n    @{ Old = \"61\u00f73=\"; New = \"46\u00f77=\" },\n    @{ Old = \"20\u00f75=\"; New = \"62\u00f72=\" },\n    @{ Old = \"84\u00f76=\"; New = \"38\u00f74=\" },\n    @{ Old = \"83\u00f75=\"; New = \"48\u00f74=\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.Old\n    $find.Replacement.Text = $r.New\n    $find.Execute($r.Old, $false, $true, $false, $false, $false, $true, 1, $false, $r.New, 2)\n}\n\nWrite-Output \"done\"\n"}
